$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row before row 242, shifting existing rows 242-259 down to 243-260
$ws.Rows.Item(242).Insert()

# Populate the newly inserted row 242 with the new data
$ws.Cells.Item(242, 1).Value = "Species name"
$ws.Cells.Item(242, 2).Value = "QUALITY"
$ws.Cells.Item(242, 3).Value = 1
$ws.Cells.Item(242, 4).Value = 1
